# Updated cryptos list on Sat Oct 21 15:30:42 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.656.56"
$ws.Range("D3").Value = "1.611.37"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").Value = "'0.994"
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'212.56"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").Value = "'0.521"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("D7").Value = "'0.994"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'28.80"
$ws.Range("E8").Value = "  +7.06%  "
$ws.Range("E9").Value = "  +3.04%  "
$ws.Range("E10").Value = "  +1.49%  "
$ws.Range("D11").Value = "'0.0908"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("D12").Value = "1.840.77"
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").Value = "1.607.73"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "'0.570"
$ws.Range("E14").Value = "  +6.21%  "
$ws.Range("E15").Value = "  +3.10%  "
$ws.Range("D16").Value = "29.687.79"
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'8.66"
$ws.Range("E17").Value = "  +13.71%  "
$ws.Range("D18").Value = "'64.64"
$ws.Range("E18").Value = "  +1.74%  "
$ws.Range("D19").Value = "'240.44"
$ws.Range("E19").Value = "  -0.98%  "
$ws.Range("E20").Value = "  +1.80%  "
$ws.Range("D21").Value = "'0.995"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("D23").Value = "'9.61"
$ws.Range("E23").Value = "  +4.81%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  +2.01%  "
$ws.Range("D25").Value = "'156.59"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Value = "'15.59"
$ws.Range("E26").Value = "  +1.82%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'6.56"
$ws.Range("E28").Value = "  +2.55%  "
$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("E31").Value = "  +2.37%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("D34").Value = "1.441.08"
$ws.Range("E34").Value = "  +1.72%  "
$ws.Range("E35").Value = "  +5.48%  "
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("E37").Value = "  +3.33%  "
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("E40").Value = "  +3.65%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.01"
$ws.Range("E41").Value = "  +2.29%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.0506"
$ws.Range("E42").Value = "  +6.69%  "
$ws.Range("E43").Value = "  +4.00%  "
$ws.Range("D44").Value = "'54.09"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").Value = "'69.79"
$ws.Range("E45").Value = "  +6.09%  "
$ws.Range("D46").Value = "'0.994"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +20.20%  "
$ws.Range("E48").Value = "  +3.12%  "
$ws.Range("D49").Value = "1.750.37"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").Value = "'87.46"
$ws.Range("E50").Value = "  +1.24%  "
$ws.Range("E51").Value = "  -0.86%  "
